$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C5").Value = "OLED CS"
$ws.Range("D7").Value = "Relay #1"
$ws.Range("G7").Value = "LED Green?"
$ws.Range("D10").Value = "DHT11, Relay #4"
$ws.Range("D8").Value = "Relay #2"
$ws.Range("D9").Value = "Relay #3"
$ws.Range("D17").Value = "DHT11?"

$ws.Columns("D").ColumnWidth = 16.33

$ws.Activate() | Out-Null
$ws.Range("D17").Select() | Out-Null
